$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp in A1
$ws.Range("A1").Value = "Datos actualizados a 6 de Octubre de 2020 a las 01:43"

# Row 4
$ws.Range("B4").Value = 7676378
$ws.Range("C4").Value = 38566
$ws.Range("D4").Value = 4888499
$ws.Range("E4").Value = 2572909
$ws.Range("G4").Value = 359
$ws.Range("H4").Value = 214970

# Row 6
$ws.Range("B6").Value = 4940499
$ws.Range("C6").Value = 25210
$ws.Range("E6").Value = 498424
$ws.Range("G6").Value = 398
$ws.Range("H6").Value = 146773

# Row 8
$ws.Range("B8").Value = 862158
$ws.Range("C8").Value = 7106
$ws.Range("D8").Value = 766300
$ws.Range("E8").Value = 69014
$ws.Range("G8").Value = 132
$ws.Range("H8").Value = 26844

# Row 37
$ws.Range("B37").Value = 115919
$ws.Range("C37").Value = 633
$ws.Range("D37").Value = 92423
$ws.Range("E37").Value = 21066
$ws.Range("G37").Value = 7
$ws.Range("H37").Value = 2430

# Row 49 -> Chequia
$ws.Range("A49").Value = "Chequia"
$ws.Range("B49").Value = 85566
$ws.Range("C49").Value = 3120
$ws.Range("D49").Value = 46636
$ws.Range("E49").Value = 38172
$ws.Range("G49").Value = 31
$ws.Range("H49").Value = 758

# Row 50 -> China
$ws.Range("A50").Value = "China"
$ws.Range("B50").Value = 85470
$ws.Range("C50").Value = 20
$ws.Range("D50").Value = 80628
$ws.Range("E50").Value = 208
$ws.Range("G50").Value = 0
$ws.Range("H50").Value = 4634

# Row 95
$ws.Range("B95").Value = 14605
$ws.Range("C95").Value = 148
$ws.Range("E95").Value = 3140

# Row 105
$ws.Range("B105").Value = 10567
$ws.Range("C105").Value = 37
$ws.Range("D105").Value = 9427
$ws.Range("E105").Value = 1106

# Row 115
$ws.Range("B115").Value = 7523
$ws.Range("C115").Value = 3
$ws.Range("D115").Value = 7204
$ws.Range("E115").Value = 157

# Row 130
$ws.Range("B130").Value = 4954
$ws.Range("C130").Value = 13
$ws.Range("D130").Value = 4755
$ws.Range("E130").Value = 93

# Row 150 -> Guinea-Bisau
$ws.Range("A150").Value = "Guinea-Bisau"
$ws.Range("B150").Value = 2385
$ws.Range("C150").Value = 23
$ws.Range("D150").Value = 1728
$ws.Range("E150").Value = 617
$ws.Range("G150").Value = 1
$ws.Range("H150").Value = 40

# Row 151 -> Principado de Andorra
$ws.Range("A151").Value = "Principado de Andorra"
$ws.Range("B151").Value = 2370
$ws.Range("C151").Value = 260
$ws.Range("D151").Value = 1615
$ws.Range("E151").Value = 702
$ws.Range("H151").Value = 53

# Row 154 -> Polinesia Francesa
$ws.Range("A154").Value = "Polinesia Francesa"
$ws.Range("B154").Value = 2228
$ws.Range("C154").Value = 264
$ws.Range("D154").Value = 1769
$ws.Range("E154").Value = 450
$ws.Range("H154").Value = 9

# Row 155 -> Belice
$ws.Range("A155").Value = "Belice"
$ws.Range("B155").Value = 2196
$ws.Range("C155").Value = 65
$ws.Range("D155").Value = 1378
$ws.Range("E155").Value = 788
$ws.Range("G155").Value = 1
$ws.Range("H155").Value = 30

# Row 156 -> Burkina Faso
$ws.Range("A156").Value = "Burkina Faso"
$ws.Range("B156").Value = 2184
$ws.Range("C156").Value = 17
$ws.Range("D156").Value = 1420
$ws.Range("E156").Value = 705
$ws.Range("H156").Value = 59

# Row 157 -> Uruguay
$ws.Range("A157").Value = "Uruguay"
$ws.Range("B157").Value = 2155
$ws.Range("C157").Value = 10
$ws.Range("D157").Value = 1862
$ws.Range("E157").Value = 245
$ws.Range("G157").Value = 0
$ws.Range("H157").Value = 48

# Row 158 -> Letonia
$ws.Range("A158").Value = "Letonia"
$ws.Range("B158").Value = 2126
$ws.Range("C158").Value = 40
$ws.Range("D158").Value = 1307
$ws.Range("E158").Value = 780
$ws.Range("H158").Value = 39

# Row 159 -> Yemen
$ws.Range("A159").Value = "Yemen"
$ws.Range("B159").Value = 2041
$ws.Range("D159").Value = 1323
$ws.Range("E159").Value = 126
$ws.Range("G159").Value = 1
$ws.Range("H159").Value = 592

# Row 184
$ws.Range("B184").Value = 387
$ws.Range("C184").Value = 2
$ws.Range("E184").Value = 20
